$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D values (historical_growth_revenue_last_5_years)
$ws.Range("D2").Value = -0.0373
$ws.Range("D3").Value = -0.0373

# Row 2 updates
$ws.Range("G2").Value = -0.04081081081081081
$ws.Range("H2").Value = -0.04081081081081081
$ws.Range("I2").Value = -0.1045045045045045
$ws.Range("J2").Value = -0.1045045045045045
$ws.Range("K2").Value = -4.18
$ws.Range("L2").Value = -0.1882882882882883
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = -0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = -0
$ws.Range("T2").ClearContents()
$ws.Range("U2").Value = 1.22
$ws.Range("V2").Value = 0.6256410256410256
$ws.Range("W2").Value = 0.2271739130434783
$ws.Range("X2").Value = 1.400551290075521
$ws.Range("Y2").Value = -1.173377377032043
$ws.Range("Z2").Value = 1.043968963084881
$ws.Range("AA2").Value = -0.1090994592052668
$ws.Range("AB2").Value = 0.111001412584427
$ws.Range("AC2").Value = -0.2201008717896939
$ws.Range("AD2").Value = 43
$ws.Range("AF2").Value = 43
$ws.Range("AG2").Value = 41.78
$ws.Range("AH2").Value = 0.9566184649610678
$ws.Range("AI2").Value = 2.028301886792453
$ws.Range("AJ2").Value = 0.9554081865995884
$ws.Range("AK2").Value = 2.091091091091091
$ws.Range("AL2").Value = 1.93
$ws.Range("AM2").Value = 1.917
$ws.Range("AN2").Value = -42.57425742574257
$ws.Range("AO2").Value = -1.202072538860104
$ws.Range("AP2").Value = -41.36633663366337
$ws.Range("AQ2").Value = -1.210224308815858

# Row 3 updates
$ws.Range("G3").Value = -0.04081081081081081
$ws.Range("H3").Value = -0.04081081081081081
$ws.Range("I3").Value = -0.1045045045045045
$ws.Range("J3").Value = -0.1045045045045045
$ws.Range("K3").Value = -4.18
$ws.Range("L3").Value = -0.1882882882882883
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = 0
$ws.Range("T3").ClearContents()
$ws.Range("U3").Value = 1.22
$ws.Range("V3").Value = 0.6256410256410256
$ws.Range("W3").Value = 0.2271739130434783
$ws.Range("X3").Value = 1.400551290075521
$ws.Range("Y3").Value = -1.173377377032043
$ws.Range("Z3").Value = 1.043968963084881
$ws.Range("AA3").Value = -0.1090994592052668
$ws.Range("AB3").Value = 0.111001412584427
$ws.Range("AC3").Value = -0.2201008717896939
$ws.Range("AD3").Value = 43
$ws.Range("AF3").Value = 43
$ws.Range("AG3").Value = 41.78
$ws.Range("AH3").Value = 0.9566184649610678
$ws.Range("AI3").Value = 2.028301886792453
$ws.Range("AJ3").Value = 0.9554081865995884
$ws.Range("AK3").Value = 2.091091091091091
$ws.Range("AL3").Value = 1.93
$ws.Range("AM3").Value = 1.917
$ws.Range("AN3").Value = -42.57425742574257
$ws.Range("AO3").Value = -1.202072538860104
$ws.Range("AP3").Value = -41.36633663366337
$ws.Range("AQ3").Value = -1.210224308815858

Write-Host "Capital structure database updated"
